# Build site at 2022-09-26 16:07:08 UTC
# - Delete row 13 (old "Docentes responsáveis:" value row), shifting rows 14+
#   up by one.
# - Update a handful of cells whose text content was re-entered by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 13 entirely; everything below shifts up one row.
$ws.Rows(13).Delete()

# Row 10 ("Objetivos:") now carries the professor's name instead of the
# long objectives paragraph.
$ws.Range("B10").Value = "4808662 - Lucrécio Fábio dos Santos"
$ws.Range("C10").Value = "4808662 - Lucrécio Fábio dos Santos"

# Row 13 ("Programa resumido:") now just says "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 ("Programa:") now holds the activation date.
$ws.Range("B15").Value = "01/01/2022"
$ws.Range("C15").Value = "01/01/2022"

# Row 18 ("Método:") now carries the professor's name.
$ws.Range("B18").Value = "4808662 - Lucrécio Fábio dos Santos"
$ws.Range("C18").Value = "4808662 - Lucrécio Fábio dos Santos"

# Row 19 ("Critério:") now holds the teaching-method text.
$ws.Range("B19").Value = "Aulas expositivas, desenvolvimento de trabalhos e exercícios em sala e fora de sala de aula, discussão de casos práticos."
$ws.Range("C19").Value = "Aulas expositivas, desenvolvimento de trabalhos e exercícios em sala e fora de sala de aula, discussão de casos práticos."

# Row 20 ("Norma de recuperação:") now holds the evaluation-criteria text.
$ws.Range("B20").Value = "Provas em sala, entrega de trabalhos e exercícios ou casos práticos elaborados fora de sala de aula."
$ws.Range("C20").Value = "Provas em sala, entrega de trabalhos e exercícios ou casos práticos elaborados fora de sala de aula."

# Row 21 ("Bibliografia:") now holds the make-up-exam rule text instead of
# the bibliography.
$ws.Range("B21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
$ws.Range("C21").Value = "Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita prova de recuperação."
